# Implementat el generador i avaluador d'Oliva (1992) que diferencia els
# valors dels batecs accentuals (0-3).
#
# Re-derives the example table (rows 3-11): the ten-syllable strings keep
# the same set of values for RC2 / RC3 / Complexitat and the accent-beat
# position lists (columns F and G), but are re-ordered/updated to reflect
# the new generator+evaluator logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that look like plain numbers ("1", "3", "7", ...) need a leading
# apostrophe so Excel stores them as text (matching the original sheet,
# where the position-list columns are text, not numeric).
$rows = @(
    @{ Row = 3;  A = "ATAAATATAT"; C = 1; D = 0; E = 1; F = "'3";    G = ""   },
    @{ Row = 4;  A = "ATATATAAAT"; C = 1; D = 0; E = 1; F = "'7";    G = ""   },
    @{ Row = 5;  A = "TAATATATAT"; C = 0; D = 1; E = 1; F = "";      G = "'1" },
    @{ Row = 6;  A = "AAATATATAT"; C = 1; D = 0; E = 1; F = "'1";    G = ""   },
    @{ Row = 7;  A = "TAAAATATAT"; C = 1; D = 1; E = 2; F = "'3";    G = "'1" },
    @{ Row = 8;  A = "AAAAATATAT"; C = 2; D = 0; E = 2; F = "1, 3";  G = ""   },
    @{ Row = 9;  A = "AAATATAAAT"; C = 2; D = 0; E = 2; F = "1, 7";  G = ""   },
    @{ Row = 10; A = "ATAAATAAAT"; C = 2; D = 0; E = 2; F = "3, 7";  G = ""   },
    @{ Row = 11; A = "TAATATAAAT"; C = 1; D = 1; E = 2; F = "'7";    G = "'1" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value2 = $r.A
    $ws.Cells.Item($r.Row, 3).Value2 = $r.C
    $ws.Cells.Item($r.Row, 4).Value2 = $r.D
    $ws.Cells.Item($r.Row, 5).Value2 = $r.E
    $ws.Cells.Item($r.Row, 6).Value  = $r.F
    $ws.Cells.Item($r.Row, 7).Value  = $r.G
}
